$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($cellRef, [string]$val) {
    $c = $ws.Range($cellRef)
    $c.NumberFormat = "@"
    $c.Value = $val
    $c.Style = "Normal"
}

Set-TextValue 'D2' '26.323.62'
Set-TextValue 'E2' '  -0.83%  '
Set-TextValue 'D3' '1.714.14'
Set-TextValue 'E3' '  -1.45%  '
Set-TextValue 'D4' '0.9970'
Set-TextValue 'E4' '  -0.24%  '
Set-TextValue 'D5' '240.47'
Set-TextValue 'E5' '  -2.65%  '
Set-TextValue 'D6' '0.9977'
Set-TextValue 'E6' '  -0.19%  '
Set-TextValue 'D7' '0.4860'
Set-TextValue 'E7' '  -0.89%  '
Set-TextValue 'D8' '0.2583'
Set-TextValue 'E8' '  -3.30%  '
Set-TextValue 'D9' '0.06161'
Set-TextValue 'E9' '  -2.32%  '
Set-TextValue 'D10' '1.711.81'
Set-TextValue 'E10' '  -1.25%  '
Set-TextValue 'D11' '0.06937'
Set-TextValue 'E11' '  -1.48%  '
Set-TextValue 'D12' '15.45'
Set-TextValue 'D13' '4.455'
Set-TextValue 'E13' '  -3.36%  '
Set-TextValue 'D14' '0.5952'
Set-TextValue 'E14' '  -2.63%  '
Set-TextValue 'D15' '76.33'
Set-TextValue 'E15' '  -1.48%  '
Set-TextValue 'D16' '0.9975'
Set-TextValue 'E16' '  -0.20%  '
Set-TextValue 'B17' 'BinanceUSD'
Set-TextValue 'C17' 'https://coinranking.com/coin/vSo2fu9iE1s0Y+binanceusd-busd'
Set-TextValue 'D17' '0.9970'
Set-TextValue 'E17' '  -0.29%  '
Set-TextValue 'B18' 'WrappedBTC'
Set-TextValue 'C18' 'https://coinranking.com/coin/x4WXHge-vvFY+wrappedbtc-wbtc'
Set-TextValue 'D18' '26.236.21'
Set-TextValue 'E18' '  -1.12%  '
Set-TextValue 'D19' '0.000007066'
Set-TextValue 'E19' '  -4.71%  '
Set-TextValue 'D20' '11.21'
Set-TextValue 'E20' '  -2.89%  '
Set-TextValue 'D21' '1.933.96'
Set-TextValue 'E21' '  -1.23%  '
Set-TextValue 'D22' '4.400'
Set-TextValue 'E22' '  -4.05%  '
Set-TextValue 'D23' '8.386'
Set-TextValue 'E23' '  -3.84%  '
Set-TextValue 'D24' '5.032'
Set-TextValue 'E24' '  -4.24%  '
Set-TextValue 'D25' '135.98'
Set-TextValue 'E25' '  -3.19%  '
Set-TextValue 'D26' '15.18'
Set-TextValue 'E26' '  -1.94%  '
Set-TextValue 'D27' '1.396'
Set-TextValue 'E27' '  -1.50%  '
Set-TextValue 'B28' 'BitcoinCash'
Set-TextValue 'C28' 'https://coinranking.com/coin/ZlZpzOJo43mIo+bitcoincash-bch'
Set-TextValue 'D28' '105.80'
Set-TextValue 'E28' '  -2.01%  '
Set-TextValue 'B29' 'LidoDAOToken'
Set-TextValue 'C29' 'https://coinranking.com/coin/Pe93bIOD2+lidodaotoken-ldo'
Set-TextValue 'D29' '1.721'
Set-TextValue 'E29' '  -2.66%  '
Set-TextValue 'D30' '3.861'
Set-TextValue 'E30' '  -4.61%  '
Set-TextValue 'D31' '0.07940'
Set-TextValue 'E31' '  -1.39%  '
Set-TextValue 'D32' '3.600'
Set-TextValue 'E32' '  -3.40%  '
Set-TextValue 'D33' '0.04422'
Set-TextValue 'E33' '  -3.76%  '
Set-TextValue 'B34' 'HuobiToken'
Set-TextValue 'C34' 'https://coinranking.com/coin/DXwP4wF9ksbBO+huobitoken-ht'
Set-TextValue 'D34' '2.600'
Set-TextValue 'E34' '  -0.35%  '
Set-TextValue 'B35' 'ARBITRUM'
Set-TextValue 'C35' 'https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb'
Set-TextValue 'D35' '0.9866'
Set-TextValue 'E35' '  -2.23%  '
Set-TextValue 'B36' 'ImmutableX'
Set-TextValue 'C36' 'https://coinranking.com/coin/Z96jIvLU7+immutablex-imx'
Set-TextValue 'D36' '0.6151'
Set-TextValue 'E36' '  -3.64%  '
Set-TextValue 'B37' 'TrustWalletToken'
Set-TextValue 'C37' 'https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt'
Set-TextValue 'D37' '0.9295'
Set-TextValue 'E37' '  +3.94%  '
Set-TextValue 'B38' 'RenderToken'
Set-TextValue 'C38' 'https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr'
Set-TextValue 'D38' '1.972'
Set-TextValue 'E38' '  -2.59%  '
Set-TextValue 'B39' 'MXToken'
Set-TextValue 'C39' 'https://coinranking.com/coin/QUC5kVAxSoB-+mxtoken-mx'
Set-TextValue 'D39' '2.368'
Set-TextValue 'E39' '  -1.65%  '
Set-TextValue 'B40' 'PaxDollar'
Set-TextValue 'C40' 'https://coinranking.com/coin/JCKLgWPAF+paxdollar-usdp'
Set-TextValue 'D40' '0.9968'
Set-TextValue 'E40' '  -0.80%  '
Set-TextValue 'B41' 'VeChain'
Set-TextValue 'C41' 'https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet'
Set-TextValue 'D41' '0.01466'
Set-TextValue 'E41' '  -2.70%  '
Set-TextValue 'B42' 'Quant'
Set-TextValue 'C42' 'https://coinranking.com/coin/bauj_21eYVwso+quant-qnt'
Set-TextValue 'D42' '99.80'
Set-TextValue 'E42' '  -2.58%  '
Set-TextValue 'B43' 'FraxShare'
Set-TextValue 'C43' 'https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs'
Set-TextValue 'D43' '5.395'
Set-TextValue 'E43' '  +0.00%  '
Set-TextValue 'B44' 'TheSandbox'
Set-TextValue 'C44' 'https://coinranking.com/coin/pxtKbG5rg+thesandbox-sand'
Set-TextValue 'D44' '0.3797'
Set-TextValue 'E44' '  -2.78%  '
Set-TextValue 'B45' 'Aptos'
Set-TextValue 'C45' 'https://coinranking.com/coin/HGYj5JCv5+aptos-apt'
Set-TextValue 'D45' '6.830'
Set-TextValue 'E45' '  -1.05%  '
Set-TextValue 'B46' 'Algorand'
Set-TextValue 'C46' 'https://coinranking.com/coin/TpHE2IShQw-sJ+algorand-algo'
Set-TextValue 'D46' '0.1146'
Set-TextValue 'E46' '  -3.44%  '
Set-TextValue 'B47' 'Cronos'
Set-TextValue 'C47' 'https://coinranking.com/coin/65PHZTpmE55b+cronos-cro'
Set-TextValue 'D47' '0.05342'
Set-TextValue 'E47' '  -1.08%  '
Set-TextValue 'B48' 'Elrond'
Set-TextValue 'C48' 'https://coinranking.com/coin/omwkOTglq+elrond-egld'
Set-TextValue 'D48' '30.41'
Set-TextValue 'E48' '  -0.64%  '
Set-TextValue 'B49' 'EnergySwap'
Set-TextValue 'C49' 'https://coinranking.com/coin/SbWqqTui-+energyswap-ens'
Set-TextValue 'D49' '7.646'
Set-TextValue 'E49' '  -1.74%  '
Set-TextValue 'B50' 'Aave'
Set-TextValue 'C50' 'https://coinranking.com/coin/ixgUfzmLR+aave-aave'
Set-TextValue 'D50' '50.91'
Set-TextValue 'E50' '  -1.93%  '
Set-TextValue 'B51' 'TrueUSD'
Set-TextValue 'C51' 'https://coinranking.com/coin/1ZZI6g5k5royD+trueusd-tusd'
Set-TextValue 'D51' '1.000'
Set-TextValue 'E51' '  -0.21%  '
